$wb = $excel.ActiveWorkbook

$sheets = @("zh-cn", "de-de")

foreach ($sheetName in $sheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the "Error Detail" column (P / column 16) to fit the new message.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    if ($sheetName -eq "zh-cn") {
        $handbackFile = "6d87f301-888a-48fb-b687-dece5f562898.467e3eaa1ab2a7aebd435e1046d86019be1888b1.zh-cn.xlf"
        $handbackDate = "2016-09-03 20:44:56"
    } else {
        $handbackFile = "6d87f301-888a-48fb-b687-dece5f562898.467e3eaa1ab2a7aebd435e1046d86019be1888b1.de-de.xlf"
        $handbackDate = "2016-09-03 20:45:08"
    }

    $targetFileName = "6d87f301-888a-48fb-b687-dece5f562898.md"
    $errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ed96a52591866db4d4d3974efa4341e19d762b1/e2e/6d87f301-888a-48fb-b687-dece5f562898.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd1f64b73e815d7c021d30ee48ab87515b6b2e85/e2e/6d87f301-888a-48fb-b687-dece5f562898.md."

    # Row 7 ("6d87f301-...") now carries a Latest Target File / Handback File /
    # Handback DateTime / Error Detail describing a stale handback.
    $ws.Range("J7").Value = $handbackFile
    $ws.Range("K7").Value = $handbackDate
    $ws.Range("P7").Value = $errorDetail

    # I7 becomes a hyperlink (like A7) pointing at the same handback markdown file.
    $ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd1f64b73e815d7c021d30ee48ab87515b6b2e85/e2e/6d87f301-888a-48fb-b687-dece5f562898.md", "", "", $targetFileName)
}
